$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column widths (character units); the stored OOXML "width" adds the
# standard ~0.833-char padding on top of this.
$ws.Columns.Item(1).ColumnWidth = 15
$ws.Columns.Item(2).ColumnWidth = 50
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 25

# Data rows (Item Number, Item Description, Qty, By air or ship, Remarks)
# for rows 2..23 -- the table was re-grouped/re-split and two rows were
# inserted (new TAC11173780/SC row, new TAC11184580/SC row), so the
# whole body is rewritten in its final order.
$data = @(
    @('TAC00070840', 'Acetate#ELDRM0001T 170x700x6.0mm(0.9kg/st)', 13.5, 'AC', 'AC 13.5, SC rest'),
    @('TAC00070840', 'Acetate#ELDRM0001T 170x700x6.0mm(0.9kg/st)', 31.9, 'SC', 'AC 13.5, SC rest'),
    @('TAC01114990', 'Acetate#ELMT1024 165x1450x4.0mm(1.3kg/st)', 1.3, 'SC', ''),
    @('TAC11173780', 'Acetate#ELSPH0025 160x690x8.0/1.17', 7.02, 'AC', ''),
    @('TAC11173780', 'Acetate#ELSPH0025 160x690x8.0/1.17', 22.3, 'SC', ''),
    @('TAC11179620', 'Acetate#ELMT2478 165x1450x6.0/1.9', 3.8, 'AC', ''),
    @('TAC11179670', 'Acetate#ELLAM2446 165x690x6.0/0.92', 22.397, 'AC', 'AC 22.397, SC rest'),
    @('TAC11179670', 'Acetate#ELLAM2446 165x690x6.0/0.92', 1.8, 'SC', 'AC 22.397, SC rest'),
    @('TAC11179680', 'Acetate#ELLAM2446 165x690x8.0/1.23', 18.312, 'AC', 'AC 18.312, SC rest'),
    @('TAC11179680', 'Acetate#ELLAM2446 165x690x8.0/1.23', 15.6, 'SC', 'AC 18.312, SC rest'),
    @('TAC11179680', 'Acetate#ELLAM2446 165x690x8.0/1.23', 1.2, '', ''),
    @('TAC11179770', 'Acetate#ELMM1703 170x1450x4.0/1.3', 6.5, 'AC', ''),
    @('TAC11181200', 'Acetate#ELSH0015 165x1450x6.0/1.9', 38, 'AC', 'AC 38, SC rest'),
    @('TAC11181200', 'Acetate#ELSH0015 165x1450x6.0/1.9', 36.3, 'SC', 'AC 38, SC rest'),
    @('TAC11181200', 'Acetate#ELSH0015 165x1450x6.0/1.9', 37.9, '', ''),
    @('TAC11181200', 'Acetate#ELSH0015 165x1450x6.0/1.9', 1.9, '', ''),
    @('TAC11181440', 'Acetate#ELMO2105 165x1450x4.0/1.3', 2.6, 'AC', 'AC 2.6, SC rest'),
    @('TAC11181440', 'Acetate#ELMO2105 165x1450x4.0/1.3', 2.6, 'SC', 'AC 2.6, SC rest'),
    @('TAC11184580', 'Acetate#ELMT2483 165x1450x8.0/2.63', 7.89, 'AC', ''),
    @('TAC11184580', 'Acetate#ELMT2483 165x1450x8.0/2.63', 30.6, 'SC', ''),
    @('TAC11187880', 'Acetate#ELMJ1080 170x1400x8.0/2.63', 27, 'SC', ''),
    @('TAF11178770', 'Bio Acetate#ELBMO3153 165x1450x4.0/1.3', 10.4, 'AC', '')

)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
